$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rename the TOC bookmark (_Toc194777141 -> _Toc194782465), keeping the
#    same range it currently covers.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_Toc194777141")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("_Toc194782465", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 2. Re-point the six bulleted requirement paragraphs from numId=1 to the
#    new list definition numId=2.
# ---------------------------------------------------------------------------
foreach ($p in @($d.Paragraphs)) {
    $pkgXml = $p.Range.WordOpenXML
    if ($pkgXml -like '*<w:numId w:val="1"/>*') {
        $startTag = $pkgXml.IndexOf("<w:p ")
        $endTag = $pkgXml.IndexOf("</w:p>") + 6
        $frag = $pkgXml.Substring($startTag, $endTag - $startTag)
        $newFrag = $frag.Replace('<w:numId w:val="1"/>', '<w:numId w:val="2"/>')
        $p.Range.InsertXML($newFrag)
    }
}

# ---------------------------------------------------------------------------
# 3. Finish the "reminder" paragraph: give it an explicit pPr/rPr (szCs 24)
#    and complete its sentence.
# ---------------------------------------------------------------------------
foreach ($p in @($d.Paragraphs)) {
    if ($p.Range.Text -like "Another integrated feature*") {
        $pkgXml = $p.Range.WordOpenXML
        $startTag = $pkgXml.IndexOf("<w:p ")
        $endTag = $pkgXml.IndexOf("</w:p>") + 6
        $frag = $pkgXml.Substring($startTag, $endTag - $startTag)

        $oldRun = '<w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>Another integrated feature should allow users to set a reminder by selecting a date. The reminder will take input on the plant type and species (e.g., succulent, tropical</w:t></w:r>'
        $newRun = '<w:pPr><w:rPr><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Another integrated feature should allow users to set a reminder by selecting a date. The reminder will take input on the plant type and species (e.g., succulent, tropical, houseplant) and notify the user of when to next water their plant. </w:t></w:r>'

        $newFrag = $frag.Replace($oldRun, $newRun)
        $p.Range.InsertXML($newFrag)
        break
    }
}

# ---------------------------------------------------------------------------
# 4. numbering.xml: add the new "lowerLetter" multilevel list definition
#    (abstractNumId 0 / numId 2) and shift the original one to
#    abstractNumId 1, keeping numId 1 pointed at it.
# ---------------------------------------------------------------------------
$newAbstractNum = @'
<w:abstractNum w:abstractNumId="0" w15:restartNumberingAfterBreak="0"><w:nsid w:val="10E374A1"/><w:multiLevelType w:val="multilevel"/><w:tmpl w:val="08DC62C0"/><w:lvl w:ilvl="0"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%1)"/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="0"/></w:tabs><w:ind w:left="720" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%2."/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="0"/></w:tabs><w:ind w:left="1440" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="2"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%3."/><w:lvlJc w:val="right"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="0"/></w:tabs><w:ind w:left="2160" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="3"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%4."/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="0"/></w:tabs><w:ind w:left="2880" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="4"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%5."/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="0"/></w:tabs><w:ind w:left="3600" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="5"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%6."/><w:lvlJc w:val="right"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="0"/></w:tabs><w:ind w:left="4320" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="6"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%7."/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="0"/></w:tabs><w:ind w:left="5040" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="7"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%8."/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="0"/></w:tabs><w:ind w:left="5760" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="8"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%9."/><w:lvlJc w:val="right"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="0"/></w:tabs><w:ind w:left="6480" w:hanging="180"/></w:pPr></w:lvl></w:abstractNum>
'@

$fullXml = $d.Content.WordOpenXML

$oldAbstractOpen = '<w:abstractNum w:abstractNumId="0" w15:restartNumberingAfterBreak="0">'
$newAbstractOpen = '<w:abstractNum w:abstractNumId="1" w15:restartNumberingAfterBreak="0">'
$fullXml = $fullXml.Replace($oldAbstractOpen, $newAbstractNum.Trim() + $newAbstractOpen)

$oldNum = '<w:num w:numId="1" w16cid:durableId="1456564468"><w:abstractNumId w:val="0"/></w:num>'
$newNum = '<w:num w:numId="1" w16cid:durableId="1456564468"><w:abstractNumId w:val="1"/></w:num><w:num w:numId="2" w16cid:durableId="958881553"><w:abstractNumId w:val="0"/></w:num>'
$fullXml = $fullXml.Replace($oldNum, $newNum)

$d.Content.InsertXML($fullXml)

Write-Host "Done"
